$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'first looks leggings',
    'first touch soccer',
    'first way leggings',
    'fitness apparel for women',
    'fitness clothes',
    'fitness clothes for women',
    'fitness compression pants men',
    'fitness gear mens',
    'fitness gear stability ball',
    'fitness gear women',
    'fitness legging',
    'fitness leggings',
    'fitness leggings for women',
    'fitness pants',
    'fitness pro',
    'fitness runners',
    'fitness tech',
    'fitness thigh',
    'fitness tights for women',
    'fitness tracker teen',
    'fitness tracker with pace',
    'fitness trackers for runners',
    'fitness travel',
    'fitness waist',
    'fitness wear',
    'fitness wear for women',
    'fitness wear for women plus size',
    'fitness wedding band men',
    'fitted capris for women',
    'fitted knee brace',
    'fitted tights for women',
    'flash compression',
    'flash costume cw',
    'flash cw costume',
    'flash cw shirt',
    'flash cw suit',
    'flash drive basketball',
    'flash jacket for boys',
    'flash leggings',
    'flash leggings men',
    'flash lift bleach',
    'flash running games',
    'flash spandex suit',
    'flash tights',
    'flash track suit',
    'flash underwear women',
    'flat abs for women',
    'flat ball frisbee',
    'flat weightlifting shoes',
    'flattering leggings for women',
    'fleece limed leggings',
    'fleece running tights women',
    'flex mt plus',
    'flight savers',
    'flight stabilization',
    'flight support',
    'floral athletic',
    'floral equipment',
    'floral jogging',
    'floral leggings',
    'floral tights',
    'floral tights for women',
    'flowered leggings',
    'flowered leggings for women',
    'flowered pants for women',
    'flying machine model',
    'flying pants',
    'fold up bike',
    'fold up bikes for women',
    'foot flexor',
    'foot leggings for women',
    'foot strap leggings',
    'football 3 4 tights',
    'football athletic tape',
    'football back flap',
    'football cold weather',
    'football compression pants nike',
    'football compression pants white',
    'football flat jacket',
    'football hand warmer red',
    'football knee brace running back',
    'football rainbow',
    'football short tights',
    'football tights blue',
    'for every mountain performance track',
    'for leggings for women',
    'for life exo',
    'forking women',
    'form fitting leggings',
    'form fitting pants for women',
    'fox riding gear set',
    'fox winter riding gear',
    'frames per second',
    'free basketball games 2017',
    'free bike race pro',
    'free leggings',
    'free run distance women',
    'free soccer games 2017',
    'freezer clothes for women',
    'freezer dress'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
